$d = $word.ActiveDocument
$normal = $d.Styles("Normal")
$f = $normal.Font
Write-Host "Bold before:" $f.Bold
$f.Bold = 9999999
Write-Host "Bold after:" $f.Bold
